$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New service-log rows appended to the bottom of Sheet1 (rows 314-323).
# Columns: A=DATE(serial) B=VECHILE REG NO C=VEHICLE BRAND D=ISSUE
#          E=STATUS F=AMOUNT G=CASH TYPE

$rows = @(
    @{ A = 44798; B = "KL01CF1995"; C = "TIGOR";      D = "RUNNING REPAIR";      E = "WORK DONE DELIVERED"; F = 550;   G = "CASH" },
    @{ A = 44798; B = "KA01MK8821"; C = "ECOSPORT";    D = "PMS";                 E = "WORK DONE" },
    @{ A = 44799; B = "KA03MM9548"; C = "I20";         D = "PMS";                 E = "WORK DONE" },
    @{ A = 44799; B = "KA01MF2461"; C = "FABIA";       D = "PMS                                      WW"; E = "WORK DONE DELIVERED"; F = 2159;  G = "CREDIT" },
    @{ A = 44799; B = "MP50BC8265"; C = "INNOVA";      D = "RUNNING REPAIR";      E = "WORK DONE DELIVERED"; F = 2500;  G = "G PAY" },
    @{ A = 44799; B = "KL05AP981";  C = "KWDI";        D = "PARTS";               E = "WORK DONE DELIVERED"; F = 2798;  G = "GPAY" },
    @{ A = 44799; B = "KA51N2602";  C = "BALENO";      D = "WIPER BLADE CHANGE";  E = "WORK DONE DELIVERED"; F = 367;   G = "P PAY" },
    @{ A = 44799; B = "KA03MS2872"; C = "BEAT";        D = "PMS";                 E = "WORK DONE DELIVERED"; F = 18582; G = "CREDIT" },
    @{ A = 44799; B = "KA03MY1985"; C = "XUV 500";     D = "TIE LAMP CHANGE";     E = "WORK DONE DELIVERED"; F = 4318;  G = "CREDIT" },
    @{ A = 44799; B = "KA01MK8821"; C = "ECOSPORT";    D = "PMS";                 E = "WORK DONE DELIVERED"; F = 4356;  G = "G PAY" }
)

$r = 314
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    if ($row.ContainsKey("F")) {
        $ws.Cells.Item($r, 6).Value = $row.F
    }
    if ($row.ContainsKey("G")) {
        $ws.Cells.Item($r, 7).Value = $row.G
    }
    $r = $r + 1
}

# Match the author's final cursor position / scrolled view.
[void]$ws.Range("F320").Select()
